# Commit: "1st changes of mifos to finflux"
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting old N/O/P -> O/P/Q.
# - Make "Repayment schedule" the active/selected sheet (was "NewLoanInput").
# - Update the selection on "Repayment schedule" to S5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at position N (14th column); this pushes the
# existing Late / heading / Outstanding columns one to the right.
$ws.Columns.Item(14).Insert()

# Match the width of the inserted column to its neighbours (closest the
# host's column-width model can represent).
$ws.Columns.Item(14).ColumnWidth = 9.8

# Make the Repayment schedule sheet the active tab and update its selection.
$ws.Activate()
$ws.Range("S5").Select() | Out-Null
